$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 88, pushing existing rows 88:109 down to 89:110.
$ws.Rows("88:88").Insert()

# Populate the newly inserted row 88 with the new weekly record.
# (Same "common" columns as the surrounding rows: Mercado ID, Mercado, Region, Codreg,
# Categoria ID, Categoria, Calidad, Clasificacion.)
$ws.Range("A88").Value = 7
$ws.Range("B88").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C88").Value = "Ñuble"
$ws.Range("D88").Value = 44889
$ws.Range("D88").NumberFormat = $ws.Range("D89").NumberFormat
$ws.Range("E88").Value = 16
$ws.Range("F88").Value = 100112021
$ws.Range("G88").Value = "Ají"
$ws.Range("H88").Value = "Americana (o)"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 60
$ws.Range("K88").Value = 16000
$ws.Range("L88").Value = 17000
$ws.Range("M88").Value = 16500
$ws.Range("N88").Value = "$/caja 15 kilos"
$ws.Range("O88").Value = "Región del Maule"
$ws.Range("P88").Value = 1100
$ws.Range("Q88").Value = 15
$ws.Range("R88").Value = "Hortaliza"
